$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B84 currently holds the text "4" (inline string); convert it to a real number 4.
$ws.Range("B84").Value = 4

# Append new row 85 with the new annotation data.
$ws.Range("A85").Value = "Ying Tang"

$ws.Range("B85").NumberFormat = "@"
$ws.Range("B85").Value = "1"

$ws.Range("C85").Value = "No technical contribution."
$ws.Range("D85").Value = "CRT"
$ws.Range("E85").Value = "THE"
$ws.Range("F85").Value = "e885cb01-c8a9-4c3f-b9a4-e5ab35292953"
$ws.Range("G85").Value = "SkwAEQbAb_annotated.xlsx"
$ws.Range("H85").Value = "No technical contribution."
